# Insert a new "statut_name" column before column C (NCTId), shifting
# columns C..L to D..M, and populate it based on the "statut_label" (column B)
# value of each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts existing columns C..L to D..M.
$ws.Range("C1").EntireColumn.Insert()

# Header for new column
$ws.Range("C1").Value2 = "statut_name"

# Determine last used row
$lastRow = $ws.Cells.Item($ws.Rows.Count, "B").End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Range("B$r").Value2
    $name = ""
    switch ($label) {
        "noir"   { $name = "pas de résultat ni de publication" }
        "rouge"  { $name = "résultat et / ou publication posté" }
        "orange" { $name = "résultat et / ou publication posté dans les 36 mois" }
        "vert"   { $name = "résultat et / ou publication posté dans les 12 mois" }
        default  { $name = "" }
    }
    $ws.Range("C$r").Value2 = $name
}
